$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-77
# from serial date 45186 (2023-09-17) to 45188 (2023-09-19)
$ws.Range("C2:C77").Value = 45188
